$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 73 & 74: the two fixtures' data got swapped (id's B, F, G, H..AC),
# while the A column (row index 71 / 72) stays put. Only write the cells
# whose value actually changes between the two rows.
# ---------------------------------------------------------------------------

# Row 73 (becomes what row 74 used to hold)
$ws.Range("B73").Value = 7646750
$ws.Range("F73").Value = "Perth Glory"
$ws.Range("G73").Value = "Wellington Phoenix"
$ws.Range("I73").Value = 4
$ws.Range("J73").Value = "A"
$ws.Range("K73").Value = 2.45
$ws.Range("L73").Value = 3.75
$ws.Range("M73").Value = 2.55
$ws.Range("N73").Value = 3.1
$ws.Range("O73").Value = 3.8
$ws.Range("P73").Value = 2.05
$ws.Range("Q73").Value = 0.25
$ws.Range("R73").Value = 2
$ws.Range("S73").Value = 1.85
$ws.Range("T73").Value = 3
$ws.Range("U73").Value = 1.925
$ws.Range("V73").Value = 1.925
$ws.Range("W73").Value = -1
$ws.Range("Y73").Value = 1.05
$ws.Range("Z73").Value = -1
$ws.Range("AA73").Value = 0.8500000000000001
$ws.Range("AB73").Value = 0.925

# Row 74 (becomes what row 73 used to hold)
$ws.Range("B74").Value = 7646749
$ws.Range("F74").Value = "Brisbane Roar"
$ws.Range("G74").Value = "Newcastle Jets"
$ws.Range("I74").Value = 2
$ws.Range("J74").Value = "H"
$ws.Range("K74").Value = 1.909
$ws.Range("L74").Value = 4
$ws.Range("M74").Value = 3.4
$ws.Range("N74").Value = 2.4
$ws.Range("O74").Value = 4
$ws.Range("P74").Value = 2.6
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = 1.83
$ws.Range("S74").Value = 2.07
$ws.Range("T74").Value = 3.25
$ws.Range("U74").Value = 1.9
$ws.Range("V74").Value = 1.95
$ws.Range("W74").Value = 1.4
$ws.Range("Y74").Value = -1
$ws.Range("Z74").Value = 0.8300000000000001
$ws.Range("AA74").Value = -1
$ws.Range("AB74").Value = 0.8999999999999999

# ---------------------------------------------------------------------------
# Rows 118 & 119: results (H/I/J) came in and odds got refreshed; two new
# settlement columns (AB/AC) appear.
# ---------------------------------------------------------------------------

$ws.Range("H118").Value = 3
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = "H"
$ws.Range("N118").Value = 1.85
$ws.Range("P118").Value = 4
$ws.Range("R118").Value = 1.85
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 2.75
$ws.Range("U118").Value = 1.825
$ws.Range("V118").Value = 2.025
$ws.Range("W118").Value = 0.8500000000000001
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.8500000000000001
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.825
$ws.Range("AC118").Value = -1

$ws.Range("H119").Value = 3
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = "H"
$ws.Range("N119").Value = 3.6
$ws.Range("O119").Value = 4
$ws.Range("P119").Value = 1.909
$ws.Range("R119").Value = 1.925
$ws.Range("S119").Value = 1.925
$ws.Range("W119").Value = 2.6
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = 0.925
$ws.Range("AA119").Value = -1
$ws.Range("AB119").Value = 1
$ws.Range("AC119").Value = -1

# ---------------------------------------------------------------------------
# Row 120: odds refreshed for an upcoming fixture (no result yet).
# ---------------------------------------------------------------------------

$ws.Range("N120").Value = 1.666
$ws.Range("O120").Value = 4.333
$ws.Range("P120").Value = 4
$ws.Range("R120").Value = 1.9
$ws.Range("S120").Value = 2
$ws.Range("T120").Value = 3.25
$ws.Range("U120").Value = 2
$ws.Range("V120").Value = 1.85

# ---------------------------------------------------------------------------
# New rows 121-125: five freshly scraped upcoming fixtures appended.
# Column A / E carry the same styles as the rest of the table, so clone
# them from row 120 before writing the new values in.
# ---------------------------------------------------------------------------

$newRows = 121, 122, 123, 124, 125
foreach ($r in $newRows) {
    $ws.Range("A120").Copy($ws.Range("A$r"))
    $ws.Range("E120").Copy($ws.Range("E$r"))
}

# Row 121
$ws.Range("A121").Value = 119
$ws.Range("B121").Value = 7126791
$ws.Range("C121").Value = "Australia ALeague"
$ws.Range("D121").Value = "Australia ALeague"
$ws.Range("E121").Value = 45360.20833333334
$ws.Range("F121").Value = "Melbourne City"
$ws.Range("G121").Value = "Wellington Phoenix"
$ws.Range("K121").Value = 2
$ws.Range("L121").Value = 3.5
$ws.Range("M121").Value = 3.6
$ws.Range("N121").Value = 1.909
$ws.Range("O121").Value = 3.8
$ws.Range("P121").Value = 3.75
$ws.Range("Q121").Value = -0.5
$ws.Range("R121").Value = 1.95
$ws.Range("S121").Value = 1.95
$ws.Range("T121").Value = 2.75
$ws.Range("U121").Value = 1.825
$ws.Range("V121").Value = 2.025
$ws.Range("W121").Value = 0
$ws.Range("X121").Value = 0
$ws.Range("Y121").Value = 0
$ws.Range("Z121").Value = 0
$ws.Range("AA121").Value = 0

# Row 122
$ws.Range("A122").Value = 120
$ws.Range("B122").Value = 7127386
$ws.Range("C122").Value = "Australia ALeague"
$ws.Range("D122").Value = "Australia ALeague"
$ws.Range("E122").Value = 45360.26041666666
$ws.Range("F122").Value = "Adelaide United"
$ws.Range("G122").Value = "Melbourne Victory"
$ws.Range("K122").Value = 2.4
$ws.Range("L122").Value = 3.4
$ws.Range("M122").Value = 2.8
$ws.Range("N122").Value = 3
$ws.Range("O122").Value = 3.5
$ws.Range("P122").Value = 2.25
$ws.Range("Q122").Value = 0.25
$ws.Range("R122").Value = 1.92
$ws.Range("S122").Value = 1.98
$ws.Range("T122").Value = 3
$ws.Range("U122").Value = 1.85
$ws.Range("V122").Value = 2
$ws.Range("W122").Value = 0
$ws.Range("X122").Value = 0
$ws.Range("Y122").Value = 0
$ws.Range("Z122").Value = 0
$ws.Range("AA122").Value = 0

# Row 123
$ws.Range("A123").Value = 121
$ws.Range("B123").Value = 7127387
$ws.Range("C123").Value = "Australia ALeague"
$ws.Range("D123").Value = "Australia ALeague"
$ws.Range("E123").Value = 45360.32291666666
$ws.Range("F123").Value = "Perth Glory"
$ws.Range("G123").Value = "Newcastle Jets"
$ws.Range("K123").Value = 1.909
$ws.Range("L123").Value = 3.75
$ws.Range("M123").Value = 3.6
$ws.Range("N123").Value = 2.15
$ws.Range("O123").Value = 3.6
$ws.Range("P123").Value = 3
$ws.Range("Q123").Value = -0.25
$ws.Range("R123").Value = 1.93
$ws.Range("S123").Value = 1.97
$ws.Range("T123").Value = 3.25
$ws.Range("U123").Value = 1.975
$ws.Range("V123").Value = 1.875
$ws.Range("W123").Value = 0
$ws.Range("X123").Value = 0
$ws.Range("Y123").Value = 0
$ws.Range("Z123").Value = 0
$ws.Range("AA123").Value = 0

# Row 124
$ws.Range("A124").Value = 122
$ws.Range("B124").Value = 7127388
$ws.Range("C124").Value = "Australia ALeague"
$ws.Range("D124").Value = "Australia ALeague"
$ws.Range("E124").Value = 45361.125
$ws.Range("F124").Value = "Sydney FC"
$ws.Range("G124").Value = "Brisbane Roar"
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 5
$ws.Range("N124").Value = 1.55
$ws.Range("O124").Value = 5
$ws.Range("P124").Value = 4.5
$ws.Range("Q124").Value = -1
$ws.Range("R124").Value = 1.9
$ws.Range("S124").Value = 2
$ws.Range("T124").Value = 3.5
$ws.Range("U124").Value = 1.975
$ws.Range("V124").Value = 1.875
$ws.Range("W124").Value = 0
$ws.Range("X124").Value = 0
$ws.Range("Y124").Value = 0
$ws.Range("Z124").Value = 0
$ws.Range("AA124").Value = 0

# Row 125
$ws.Range("A125").Value = 123
$ws.Range("B125").Value = 7128012
$ws.Range("C125").Value = "Australia ALeague"
$ws.Range("D125").Value = "Australia ALeague"
$ws.Range("E125").Value = 45361.125
$ws.Range("F125").Value = "Macarthur FC"
$ws.Range("G125").Value = "Central Coast Mariners"
$ws.Range("K125").Value = 2.4
$ws.Range("L125").Value = 3.5
$ws.Range("M125").Value = 2.75
$ws.Range("N125").Value = 3.1
$ws.Range("O125").Value = 3.5
$ws.Range("P125").Value = 2.2
$ws.Range("Q125").Value = 0.25
$ws.Range("R125").Value = 1.95
$ws.Range("S125").Value = 1.95
$ws.Range("T125").Value = 3
$ws.Range("U125").Value = 1.975
$ws.Range("V125").Value = 1.875
$ws.Range("W125").Value = 0
$ws.Range("X125").Value = 0
$ws.Range("Y125").Value = 0
$ws.Range("Z125").Value = 0
$ws.Range("AA125").Value = 0
